$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E1').Value = 'MES 1'
$ws.Range('C3').Value = 446858039.41
$ws.Range('D3').Value = 444564133.25
$ws.Range('E3').Value = 12025959.48
$ws.Range('F3').Value = 12443804.42
$ws.Range('G3').Value = 458883998.89
$ws.Range('H3').Value = 457007937.67
$ws.Range('I3').Value = 1876061.22
$ws.Range('C4').Value = 178898149.91
$ws.Range('D4').Value = 174836588.82
$ws.Range('E4').Value = 3946338
$ws.Range('F4').Value = 3278694.96
$ws.Range('G4').Value = 182844487.91
$ws.Range('H4').Value = 178115283.78
$ws.Range('I4').Value = 4729204.13
$ws.Range('C5').Value = 3315868.6
$ws.Range('D5').Value = 3030695.5
$ws.Range('E5').Value = 72281.92
$ws.Range('F5').Value = 77629.07000000001
$ws.Range('G5').Value = 3388150.52
$ws.Range('H5').Value = 3108324.57
$ws.Range('I5').Value = 279825.95
$ws.Range('D6').Value = 3392032.96
$ws.Range('F6').Value = 493.64
$ws.Range('H6').Value = 3392526.6
$ws.Range('I6').Value = 1814834.92
$ws.Range('C7').Value = 8865840.85
$ws.Range('D7').Value = 9809430.859999999
$ws.Range('E7').Value = 124426.08
$ws.Range('F7').Value = 166828.64
$ws.Range('G7').Value = 8990266.93
$ws.Range('H7').Value = 9976259.5
$ws.Range('J7').Value = 985992.5699999999
$ws.Range('C8').Value = 38965131.09
$ws.Range('E8').Value = 91664.23
$ws.Range('G8').Value = 39056795.32
$ws.Range('I8').Value = 27929059.91
$ws.Range('C9').Value = 44084628.02
$ws.Range('D9').Value = 39554245.57
$ws.Range('E9').Value = 712628.28
$ws.Range('F9').Value = 668260.53
$ws.Range('G9').Value = 44797256.3
$ws.Range('H9').Value = 40222506.1
$ws.Range('I9').Value = 4574750.2
$ws.Range('D10').Value = 33218458.2
$ws.Range('F10').Value = 666062.79
$ws.Range('H10').Value = 33884520.99
$ws.Range('J10').Value = 23585698.35
$ws.Range('C11').Value = 26927488.81
$ws.Range('D11').Value = 28418236.51
$ws.Range('E11').Value = 906514.45
$ws.Range('F11').Value = 804710.7
$ws.Range('G11').Value = 27834003.26
$ws.Range('H11').Value = 29222947.21
$ws.Range('J11').Value = 1388943.95
$ws.Range('C12').Value = 21194001.89
$ws.Range('D12').Value = 20218012.67
$ws.Range('E12').Value = 501337.04
$ws.Range('F12').Value = 562238.03
$ws.Range('G12').Value = 21695338.93
$ws.Range('H12').Value = 20780250.7
$ws.Range('I12').Value = 915088.23
$ws.Range('C13').Value = 102812310.62
$ws.Range('D13').Value = 103576692.79
$ws.Range('E13').Value = 2099171.64
$ws.Range('F13').Value = 2283237.8
$ws.Range('G13').Value = 104911482.26
$ws.Range('H13').Value = 105859930.59
$ws.Range('J13').Value = 948448.33
$ws.Range('C14').Value = 40628716.89
$ws.Range('D14').Value = 57808813.4
$ws.Range('E14').Value = 549820.03
$ws.Range('F14').Value = 283513.33
$ws.Range('G14').Value = 41178536.92
$ws.Range('H14').Value = 58092326.73
$ws.Range('J14').Value = 16913789.81
$ws.Range('C15').Value = 1782949.13
$ws.Range('D15').Value = 1795029.64
$ws.Range('E15').Value = 75478.77
$ws.Range('F15').Value = 38912.31
$ws.Range('G15').Value = 1858427.9
$ws.Range('H15').Value = 1833941.95
$ws.Range('I15').Value = 24485.95
$ws.Range('C16').Value = 0
$ws.Range('D16').Value = 0
$ws.Range('E16').Value = 1020741.64
$ws.Range('F16').Value = 20269.33
$ws.Range('G16').Value = 1020741.64
$ws.Range('H16').Value = 20269.33
$ws.Range('I16').Value = 1000472.31
$ws.Range('C17').Value = 0
$ws.Range('D17').Value = 0
$ws.Range('E17').Value = 166828.64
$ws.Range('F17').Value = 124426.08
$ws.Range('G17').Value = 166828.64
$ws.Range('H17').Value = 124426.08
$ws.Range('I17').Value = 42402.56
$ws.Range('J17').Value = 0
$ws.Range('C18').Value = 0
$ws.Range('E18').Value = 676294.02
$ws.Range('G18').Value = 676294.02
$ws.Range('I18').Value = 676294.02
$ws.Range('C19').Value = 0
$ws.Range('D19').Value = 0
$ws.Range('E19').Value = 702241.87
$ws.Range('F19').Value = 133163.83
$ws.Range('G19').Value = 702241.87
$ws.Range('I19').Value = 569078.04
$ws.Range('C20').Value = 0
$ws.Range('E20').Value = 9386.190000000001
$ws.Range('G20').Value = 9386.190000000001
$ws.Range('I20').Value = 9386.190000000001
$ws.Range('C21').Value = 0
$ws.Range('E21').Value = 96761.45
$ws.Range('G21').Value = 96761.45
$ws.Range('I21').Value = 96761.45
$ws.Range('C22').Value = 0
$ws.Range('D22').Value = 0
$ws.Range('E22').Value = 321748.82
$ws.Range('F22').Value = 10561
$ws.Range('G22').Value = 321748.82
$ws.Range('I22').Value = 311187.82
$ws.Range('C23').Value = 0
$ws.Range('E23').Value = 783437.09
$ws.Range('G23').Value = 783437.09
$ws.Range('I23').Value = 783437.09
$ws.Range('D24').Value = 0
$ws.Range('F24').Value = 3303231.16
$ws.Range('H24').Value = 3303231.16
$ws.Range('J24').Value = 3303231.16
$ws.Range('D25').Value = 0
$ws.Range('F25').Value = 14414.08
$ws.Range('H25').Value = 14414.08
$ws.Range('J25').Value = 14414.08
$ws.Range('D26').Value = 0
$ws.Range('F26').Value = 0
$ws.Range('H26').Value = 0
$ws.Range('J26').Value = 0
$ws.Range('D27').Value = 0
$ws.Range('F27').Value = 2607.94
$ws.Range('H27').Value = 2607.94
$ws.Range('J27').Value = 2607.94
$ws.Range('C28').Value = 0
$ws.Range('D28').Value = 0
$ws.Range('E28').Value = 143724.83
$ws.Range('F28').Value = 3632744.31
$ws.Range('H28').Value = 3632744.31
$ws.Range('J28').Value = 3489019.48
$ws.Range('C29').Value = 0
$ws.Range('D29').Value = 0
$ws.Range('E29').Value = 3010648.52
$ws.Range('F29').Value = 133163.83
$ws.Range('G29').Value = 3010648.52
$ws.Range('I29').Value = 2877484.69
$ws.Range('C30').Value = 0
$ws.Range('E30').Value = 142006.19
$ws.Range('G30').Value = 142006.19
$ws.Range('I30').Value = 142006.19
$ws.Range('C31').Value = 0
$ws.Range('E31').Value = 158340.78
$ws.Range('G31').Value = 158340.78
$ws.Range('I31').Value = 158340.78
$ws.Range('C32').Value = 0
$ws.Range('D32').Value = 0
$ws.Range('E32').Value = 321748.82
$ws.Range('F32').Value = 10561
$ws.Range('G32').Value = 321748.82
$ws.Range('I32').Value = 311187.82
$ws.Range('G33').Value = 958498838.16
$ws.Range('H33').Value = 960009634.36
$ws.Range('I33').Value = 49121349.47
$ws.Range('J33').Value = 50632145.67
